$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the "time_taken" timestamps on the existing "data" sheet -------
$dataSheet.Range("F2").Value  = "2021-10-05 14:21:45.502454"
$dataSheet.Range("F3").Value  = "2021-10-05 14:21:45.502462"
$dataSheet.Range("F4").Value  = "2021-10-05 14:21:45.502465"
$dataSheet.Range("F5").Value  = "2021-10-05 14:21:45.502468"
$dataSheet.Range("F6").Value  = "2021-10-05 14:21:45.502470"
$dataSheet.Range("F7").Value  = "2021-10-05 14:21:45.502473"
$dataSheet.Range("F8").Value  = "2021-10-05 14:21:45.502476"
$dataSheet.Range("F9").Value  = "2021-10-05 14:21:45.502478"
$dataSheet.Range("F10").Value = "2021-10-05 14:21:45.502481"
$dataSheet.Range("F11").Value = "2021-10-05 14:21:45.502484"
$dataSheet.Range("F12").Value = "2021-10-05 14:21:45.502486"
$dataSheet.Range("F13").Value = "2021-10-05 14:21:45.502489"
$dataSheet.Range("F14").Value = "2021-10-05 14:21:45.502491"
$dataSheet.Range("F15").Value = "2021-10-05 14:21:45.502494"
$dataSheet.Range("F16").Value = "2021-10-05 14:21:45.502496"
$dataSheet.Range("F17").Value = "2021-10-05 14:21:45.502499"
$dataSheet.Range("F18").Value = "2021-10-05 14:21:45.502501"
$dataSheet.Range("F19").Value = "2021-10-05 14:21:45.502504"
$dataSheet.Range("F20").Value = "2021-10-05 14:21:45.502506"
$dataSheet.Range("F21").Value = "2021-10-05 14:21:45.502509"
$dataSheet.Range("F22").Value = "2021-10-05 14:21:45.502511"
$dataSheet.Range("F23").Value = "2021-10-05 14:21:45.502514"
$dataSheet.Range("F24").Value = "2021-10-05 14:21:45.502516"
$dataSheet.Range("F25").Value = "2021-10-05 14:21:45.502519"
$dataSheet.Range("F26").Value = "2021-10-05 14:21:45.502521"
$dataSheet.Range("F27").Value = "2021-10-05 14:21:45.502524"
$dataSheet.Range("F28").Value = "2021-10-05 14:21:45.502527"
$dataSheet.Range("F29").Value = "2021-10-05 14:21:45.502529"
$dataSheet.Range("F30").Value = "2021-10-05 14:21:45.502531"
$dataSheet.Range("F31").Value = "2021-10-05 14:21:45.502534"

# --- Add a new "metadata" tab after "data" ---------------------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws.Name = "metadata"

# Reuse the bold/centered/bordered header style already used by "data"!B1:F1
# by copy/pasting formats only, instead of inventing new style entries.
$dataSheet.Range("B1:F1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# "data"!A2 carries the same style (bold/centered/bordered) -- reuse it for
# the row-index column of "metadata" too.
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Neurofibromatosis Type 1"
$ws.Range("C2").Value = 255

# Force "1.26" to be stored as text (not a number) like the source file,
# then drop the number-format style so the cell ends up unstyled.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.26"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "2019-06-20T15:13:27.332971Z"
$ws.Range("F2").Value = "2021-10-05 14:21:45.498752"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/255/?format=json"
